# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F4").Value  = 269
$wsExhibition.Range("F6").Value  = 43
$wsExhibition.Range("F7").Value  = 265
$wsExhibition.Range("F8").Value  = 213
$wsExhibition.Range("F9").Value  = 1972
$wsExhibition.Range("F10").Value = 348
$wsExhibition.Range("F11").Value = 4653
$wsExhibition.Range("F12").Value = 81
$wsExhibition.Range("F13").Value = 327

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value  = 269
$wsAll.Range("F8").Value  = 43
$wsAll.Range("F9").Value  = 265
$wsAll.Range("F10").Value = 213
$wsAll.Range("F13").Value = 1972
$wsAll.Range("F14").Value = 348
$wsAll.Range("F15").Value = 4653
$wsAll.Range("F16").Value = 81
$wsAll.Range("F17").Value = 327
